# AI Process Log - log a new entry (Entry #4) documenting getting the AI
# to fix the booking functionality, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 (Entry #4) currently has only the entry number filled in (A20=4).
# Copy the date-cell formatting from the row above (B19, which carries the
# "date" number format / border / alignment) onto B20 before writing the
# new date value, so the new row matches the look of rows 17-19.
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null

$ws.Range("B20").Value = 46057
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "Github copilot"
$ws.Range("E20").Value = "Got the ai to fix the booking"
$ws.Range("F20").Value = "Ai gave back the booking working completely fine now"
$ws.Range("G20").Value = "Booking worked completely fine, all of the functionality is fine on the page"
$ws.Range("H20").Value = "The ai is starting to make mistakes with its syntax, not sure if it’s the amount its trying to work it, the ai is forgetting to add { and ( to the end of its code causing errors I'd have to manually fix"
$ws.Range("I20").Value = "Did the same tests as last time and went through the code myself to find what the ai was struggling with. Tried multiple times to get the ai to fix it itself however it just wasn’t able to"

# Move the active selection to where the user left off editing.
$ws.Range("J20").Select() | Out-Null
